# Commit: Tue, May 05, 2020 10:06:48 AM
#
# 1) Swap the table style used by the 2x4 table on slide 16 from the
#    default "No Style, No Grid" style to the "Medium Style 2 - Accent 1"
#    style (GUID F8EA31E5-EFA8-465B-BCA5-728718BB73F3).
#
# 2) The presentation's two embedded themes (the slide-master theme and
#    the notes-master theme) had their roles swapped -- the slide master
#    now uses the colours that used to belong to the "Office Theme" and
#    vice-versa. The PowerPoint object model only exposes per-slide
#    theme *colours* for editing (ThemeColorScheme), so recolour the
#    slide-master theme's 12 theme colours to the "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) table style -------------------------------------------------
$s16 = $p.Slides.Item(16)
$tableShape = $s16.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{F8EA31E5-EFA8-465B-BCA5-728718BB73F3}", $true)

# --- 2) theme colours -------------------------------------------------
# Any slide can be used to reach the deck's (slide-master) theme colour
# scheme -- use slide 1.
$themeColors = $p.Slides.Item(1).ThemeColorScheme

$themeColors.Colors(1).RGB  = 0        # dk1      000000
$themeColors.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$themeColors.Colors(3).RGB  = 6968388  # dk2      44546A
$themeColors.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$themeColors.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$themeColors.Colors(6).RGB  = 3243501  # accent2  ED7D31
$themeColors.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$themeColors.Colors(8).RGB  = 49407    # accent4  FFC000
$themeColors.Colors(9).RGB  = 12874308 # accent5  4472C4
$themeColors.Colors(10).RGB = 4697456  # accent6  70AD47
$themeColors.Colors(11).RGB = 12673797 # hlink    0563C1
$themeColors.Colors(12).RGB = 7491477  # folHlink 954F72
